$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'29.193.26"
$ws.Range("E2").Value = "  -0.50%  "

# Row 3
$ws.Range("D3").Value = "'1.829.53"
$ws.Range("E3").Value = "  -0.69%  "

# Row 4
$ws.Range("D4").Value = "'0.9995"
$ws.Range("E4").Value = "  +0.13%  "

# Row 5
$ws.Range("D5").Value = "'237.69"
$ws.Range("E5").Value = "  -0.98%  "

# Row 6
$ws.Range("D6").Value = "'0.6052"
$ws.Range("E6").Value = "  -3.81%  "

# Row 7
$ws.Range("E7").Value = "  +0.10%  "

# Row 8
$ws.Range("D8").Value = "'0.07104"
$ws.Range("E8").Value = "  -4.07%  "

# Row 9
$ws.Range("D9").Value = "'0.2836"
$ws.Range("E9").Value = "  -2.28%  "

# Row 10
$ws.Range("D10").Value = "'24.09"
$ws.Range("E10").Value = "  -3.08%  "

# Row 11
$ws.Range("D11").Value = "'0.07649"
$ws.Range("E11").Value = "  -1.13%  "

# Row 12
$ws.Range("D12").Value = "'1.828.76"
$ws.Range("E12").Value = "  -0.70%  "

# Row 13
$ws.Range("D13").Value = "'4.799"
$ws.Range("E13").Value = "  -3.71%  "

# Row 14
$ws.Range("D14").Value = "'0.6419"
$ws.Range("E14").Value = "  -5.48%  "

# Row 15
$ws.Range("D15").Value = "'0.000009964"
$ws.Range("E15").Value = "  -2.77%  "

# Row 16
$ws.Range("D16").Value = "'2.081.18"
$ws.Range("E16").Value = "  -0.43%  "

# Row 17
$ws.Range("D17").Value = "'79.65"
$ws.Range("E17").Value = "  -2.86%  "

# Row 18
$ws.Range("D18").Value = "'6.009"
$ws.Range("E18").Value = "  -3.90%  "

# Row 19
$ws.Range("D19").Value = "'29.158.72"
$ws.Range("E19").Value = "  -0.57%  "

# Row 20
$ws.Range("D20").Value = "'230.85"
$ws.Range("E20").Value = "  +0.70%  "

# Row 21
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").Value = "'1.000"
$ws.Range("E21").Value = "  +0.09%  "

# Row 22
$ws.Range("B22").Value = "Avalanche"
$ws.Range("C22").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D22").Value = "'11.79"
$ws.Range("E22").Value = "  -4.35%  "

# Row 23
$ws.Range("D23").Value = "'7.027"
$ws.Range("E23").Value = "  -5.60%  "

# Row 24
$ws.Range("D24").Value = "'0.9987"
$ws.Range("E24").Value = "  -0.16%  "

# Row 25
$ws.Range("D25").Value = "'155.50"
$ws.Range("E25").Value = "  -1.78%  "

# Row 26
$ws.Range("D26").Value = "'8.055"
$ws.Range("E26").Value = "  -4.97%  "

# Row 27
$ws.Range("E27").Value = "  -4.74%  "

# Row 28
$ws.Range("E28").Value = "  -4.19%  "

# Row 29
$ws.Range("D29").Value = "'0.06890"
$ws.Range("E29").Value = "  +5.56%  "

# Row 30
$ws.Range("D30").Value = "'1.464"
$ws.Range("E30").Value = "  +1.07%  "

# Row 31
$ws.Range("D31").Value = "'1.459"
$ws.Range("E31").Value = "  -1.95%  "

# Row 32
$ws.Range("D32").Value = "'3.835"
$ws.Range("E32").Value = "  -5.64%  "

# Row 33
$ws.Range("D33").Value = "'3.806"
$ws.Range("E33").Value = "  -6.51%  "

# Row 34
$ws.Range("D34").Value = "'1.138"
$ws.Range("E34").Value = "  -0.12%  "

# Row 35
$ws.Range("D35").Value = "'1.727"
$ws.Range("E35").Value = "  -6.06%  "

# Row 36
$ws.Range("D36").Value = "'0.6598"
$ws.Range("E36").Value = "  -5.37%  "

# Row 37
$ws.Range("D37").Value = "'2.544"
$ws.Range("E37").Value = "  -1.08%  "

# Row 38
$ws.Range("D38").Value = "'1.233.56"
$ws.Range("E38").Value = "  -0.50%  "

# Row 39
$ws.Range("D39").Value = "'2.756"
$ws.Range("E39").Value = "  -2.16%  "

# Row 40
$ws.Range("D40").Value = "'0.01767"
$ws.Range("E40").Value = "  -4.89%  "

# Row 41
$ws.Range("D41").Value = "'6.580"
$ws.Range("E41").Value = "  -3.08%  "

# Row 42
$ws.Range("D42").Value = "'0.9326"
$ws.Range("E42").Value = "  -0.23%  "

# Row 43
$ws.Range("D43").Value = "'1.002"
$ws.Range("E43").Value = "  +0.28%  "

# Row 44
$ws.Range("D44").Value = "'1.990.26"
$ws.Range("E44").Value = "  -0.20%  "

# Row 45
$ws.Range("D45").Value = "'100.17"
$ws.Range("E45").Value = "  -0.53%  "

# Row 46
$ws.Range("D46").Value = "'63.42"
$ws.Range("E46").Value = "  -3.28%  "

# Row 47
$ws.Range("D47").Value = "'0.00000000118"
$ws.Range("E47").Value = "  -1.03%  "

# Row 48
$ws.Range("D48").Value = "'1.633"
$ws.Range("E48").Value = "  -4.72%  "

# Row 49
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'8.553"
$ws.Range("E49").Value = "  -4.92%  "

# Row 50
$ws.Range("B50").Value = "Aptos"
$ws.Range("C50").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D50").Value = "'6.549"
$ws.Range("E50").Value = "  -7.19%  "

# Row 51
$ws.Range("D51").Value = "'0.05589"
$ws.Range("E51").Value = "  -1.38%  "
